$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 / H3 - rejection message text simplified (hyphenated line-breaks removed)
$ws.Range("H3").Value = "Μήνυμα απόρριψης αιτήματος Πελάτη: `"Πρέπει να επιλέξετε το υποχρεωτικό πεδίο`""

# Row 13 / H13 - successful rejection message, merged onto two lines instead of three
$ws.Range("H13").Value = "Επιτυχής Απόρριψη αντικειμένου`nκαι επιστροφή στον κατάλογο του Καταστήματος. "

# Row 25 / H25 - error message now ends with a period
$ws.Range("H25").Value = "Μήνυμα μη επιτρεπτού ορίου`nχαρακτήρων: `"Το σχόλιο είναι πολύ μεγάλο`"."

# Row 16 / G16 - product data now also states the comment length used in this case
$ws.Range("G16").Value = "Προϊόν:UI`nΣχόλιο 150 χαρακτήρων"
$ws.Range("G16").WrapText = $true

# Row 19 / G19
$ws.Range("G19").Value = "Προϊόν:UI`nΑ"
$ws.Range("G19").WrapText = $true

# Row 22 / G22
$ws.Range("G22").Value = "Προϊόν:UI`nΣχόλιο 2,147,483,647`nχαρακτήρων "
$ws.Range("G22").WrapText = $true

# Row 25 / G25
$ws.Range("G25").Value = "Προϊόν:UI`nΣχόλιο 2,147,483,648`nχαρακτήρων"
$ws.Range("G25").WrapText = $true

# Move the saved selection to match the author's last position
$ws.Range("G30").Select()
